# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the existing header formatting (bold, border,
# centered) from an existing header cell so the new headers match style s="1".
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2 through 40): every row gets the same season record.
$lastRow = 40
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 93   # column AD = Wins
    $ws.Cells.Item($r, 31).Value = 69   # column AE = Losses
    $ws.Cells.Item($r, 32).Value = 0    # column AF = Ties
}
